# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-05-28 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-29 Wednesday", 2)

# Update the division problems inside the table. The table has 20 rows total:
# every 4th row (1, 5, 9, 13, 17 in 1-based indexing) holds the 5 visible
# problems, the rest are blank spacer rows. We address each cell directly by
# (row, column) so that values which coincidentally collide with other old/new
# values elsewhere in the table are not double-replaced.
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("94÷6=", "57÷2=", "32÷7=", "33÷8=", "46÷7=")
    5  = @("91÷5=", "88÷3=", "63÷5=", "50÷4=", "51÷2=")
    9  = @("61÷2=", "21÷2=", "93÷8=", "70÷4=", "86÷5=")
    13 = @("46÷3=", "62÷4=", "54÷8=", "84÷5=", "17÷2=")
    17 = @("60÷6=", "66÷6=", "15÷7=", "19÷3=", "33÷7=")
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
